$d = $word.ActiveDocument

$replacements = @(
    @{old="166×4=664"; new="627×3=1881"},
    @{old="188×3=564"; new="125×3=375"},
    @{old="972×3=2916"; new="526×8=4208"},
    @{old="236×8=1888"; new="291×2=582"},
    @{old="727×8=5816"; new="335×4=1340"},
    @{old="949×6=5694"; new="304×7=2128"},
    @{old="142×5=710"; new="363×3=1089"},
    @{old="705×9=6345"; new="607×2=1214"},
    @{old="187×7=1309"; new="971×5=4855"},
    @{old="420×7=2940"; new="687×3=2061"},
    @{old="751×3=2253"; new="364×3=1092"},
    @{old="462×6=2772"; new="448×8=3584"},
    @{old="907×4=3628"; new="379×8=3032"},
    @{old="222×3=666"; new="499×7=3493"},
    @{old="707×4=2828"; new="744×7=5208"},
    @{old="273×5=1365"; new="963×7=6741"},
    @{old="205×9=1845"; new="852×2=1704"},
    @{old="976×5=4880"; new="154×5=770"},
    @{old="540×6=3240"; new="246×3=738"},
    @{old="268×3=804"; new="949×2=1898"},
    @{old="390×9=3510"; new="408×7=2856"},
    @{old="524×6=3144"; new="390×3=1170"},
    @{old="220×7=1540"; new="669×2=1338"},
    @{old="847×8=6776"; new="732×2=1464"},
    @{old="729×5=3645"; new="854×9=7686"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
